# Update cryptos list figures (price + 1h volume change) per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "37.853.55"
$ws.Cells.Item(2, 5).Value = "  +1.32%  "

$ws.Cells.Item(3, 4).Value = "2.086.10"
$ws.Cells.Item(3, 5).Value = "  +1.04%  "

$ws.Cells.Item(4, 5).Value = "  -0.09%  "

$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "232.75"
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -0.55%  "

$ws.Cells.Item(6, 5).Value = "  -0.15%  "

$ws.Cells.Item(7, 5).Value = "  -0.09%  "

$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = "57.41"
$cell.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +0.97%  "

$ws.Cells.Item(9, 5).Value = "  +1.61%  "

$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0779"
$cell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +2.16%  "

$ws.Cells.Item(11, 5).Value = "  +2.94%  "

$ws.Cells.Item(12, 4).Value = "2.382.22"
$ws.Cells.Item(12, 5).Value = "  +0.47%  "

$ws.Cells.Item(13, 5).Value = "  -1.49%  "

$ws.Cells.Item(14, 5).Value = "  +1.95%  "

$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.760"
$cell.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  -2.12%  "

$ws.Cells.Item(16, 5).Value = "  +1.82%  "

$ws.Cells.Item(17, 4).Value = "2.080.66"
$ws.Cells.Item(17, 5).Value = "  +0.64%  "

$ws.Cells.Item(18, 4).Value = "37.765.33"
$ws.Cells.Item(18, 5).Value = "  +1.19%  "

$ws.Cells.Item(19, 5).Value = "  -1.65%  "

$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = "70.79"
$cell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +1.89%  "

$ws.Cells.Item(21, 5).Value = "  +1.29%  "

$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = "228.11"
$cell.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +0.79%  "

$ws.Cells.Item(23, 5).Value = "  -0.06%  "

$ws.Cells.Item(24, 5).Value = "  -1.78%  "

$ws.Cells.Item(25, 5).Value = "  -0.85%  "

$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = "170.64"
$cell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +1.94%  "

$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.138"
$cell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +8.91%  "

$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = "8.93"
$cell.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +1.99%  "

$ws.Cells.Item(29, 5).Value = "  +0.32%  "

$ws.Cells.Item(30, 5).Value = "  +2.03%  "

$ws.Cells.Item(31, 5).Value = "  +1.09%  "

$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.62"
$cell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +3.67%  "

$ws.Cells.Item(33, 5).Value = "  +1.54%  "

$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.59"
$cell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +0.70%  "

$ws.Cells.Item(35, 5).Value = "  +0.70%  "

$ws.Cells.Item(36, 5).Value = "  +3.47%  "

$ws.Cells.Item(37, 5).Value = "  +4.56%  "

$ws.Cells.Item(38, 5).Value = "  -0.10%  "

$ws.Cells.Item(39, 5).Value = "  -4.41%  "

$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0995"
$cell.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +5.82%  "

$ws.Cells.Item(41, 5).Value = "  -0.78%  "

$ws.Cells.Item(42, 5).Value = "  +1.04%  "

$ws.Cells.Item(43, 5).Value = "  +0.25%  "

$ws.Cells.Item(44, 4).Value = "1.450.69"
$ws.Cells.Item(44, 5).Value = "  -0.92%  "

$ws.Cells.Item(45, 5).Value = "  -1.06%  "

$ws.Cells.Item(46, 5).Value = "  +3.08%  "

$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.05"
$cell.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -6.61%  "

$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = "15.63"
$cell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +4.01%  "

$ws.Cells.Item(49, 5).Value = "  +3.25%  "

$ws.Cells.Item(50, 5).Value = "  +1.45%  "

$ws.Cells.Item(51, 4).Value = "2.277.49"
$ws.Cells.Item(51, 5).Value = "  +0.75%  "
